$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a "Measured"/"Maximum" label row above each of the four
# sub-tables (rows 9, 18, 26, 34 -> new rows 8, 17, 25, 33).
$rows = @(8, 17, 25, 33)
foreach ($r in $rows) {
    $ws.Range("C$r").Value = "Measured"
    $ws.Range("I$r").Value = "Maximum"
}
